$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.170.46"
$ws.Range("E2").Value = "  -3.59%  "
$ws.Range("D3").Value = "3.151.39"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.11"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.44%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.151.32"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("E9").Value = "  -3.33%  "
$ws.Range("E10").Value = "  -6.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.475"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.08%  "
$ws.Range("E13").Value = "  -5.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.06"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.06%  "
$ws.Range("D15").Value = "3.670.07"
$ws.Range("E15").Value = "  -2.62%  "
$ws.Range("D16").Value = "64.253.93"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "3.150.54"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("E23").Value = "  -3.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.80"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.10"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.43"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.52%  "
$ws.Range("E29").Value = "  -6.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.119"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -28.75%  "
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.95%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.20"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.10%  "
$ws.Range("E35").Value = "  -5.39%  "
$ws.Range("E36").Value = "  -5.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.12"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").Value = "0.0₃0724"
$ws.Range("E38").Value = "  -8.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "453.13"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -7.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.51%  "
$ws.Range("E41").Value = "  -5.63%  "
$ws.Range("E42").Value = "  -6.76%  "
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").Value = "2.846.89"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("E45").Value = "  -8.19%  "
$ws.Range("E46").Value = "  -7.86%  "
$ws.Range("E47").Value = "  -6.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.95%  "
